# Generate Report for Handback
#
# The "96ea4119-27f9-41b6-9d33-b6b848f64680" handback row (row 3) on both the
# zh-cn and de-de language sheets gets fresh handoff/handback timestamps now
# that the report has been (re)generated for this handback.
#
#   zh-cn!E3 (Correspond Handoff Datetime)  : 2016-03-21 12:39:44 -> 2016-03-21 12:40:51
#   zh-cn!H3 (Correspond Handback DateTime) : 2016-03-21 12:40:15 -> 2016-03-21 12:41:14
#   de-de!E3 (Correspond Handoff Datetime)  : 2016-03-21 12:39:48 -> 2016-03-21 12:40:54
#   de-de!H3 (Correspond Handback DateTime) : 2016-03-21 12:40:23 -> 2016-03-21 12:41:20

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-21 12:40:51"
$wsZhCn.Range("H3").Value = "2016-03-21 12:41:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-21 12:40:54"
$wsDeDe.Range("H3").Value = "2016-03-21 12:41:20"
